$p = $ppt.ActivePresentation
$s = $p.Slides.Item(19)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

# --- Edit 1: append remark to the " otherwise." run (paragraph 1) ---
$para1 = $tr.Paragraphs(1, 1)
$run = $para1.Runs($para1.Runs().Count, 1)
$run.Text = " otherwise. Whitespace is separator."

# --- Edit 2: merge " on " + "Windows)." into a single run (paragraph 4) ---
$para4 = $tr.Paragraphs(4, 1)
$full = $para4.Text
$needle = " on Windows)."
$idx = $full.IndexOf($needle)
$startPos = $idx + 1
$sub = $para4.Characters($startPos, $needle.Length)
$sub.Text = $needle
